$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the JLCPCB part number for the through-hole audio jack connectors
# (Cassette in, Audio out, Printer, Cassette out) from C2884998 to C2939180,
# per commit: "Update to Use throughhole audio jacks (C2939180)"
$ws.Range("D22").Value = "C2939180"
$ws.Range("D34").Value = "C2939180"
$ws.Range("D47").Value = "C2939180"
$ws.Range("D51").Value = "C2939180"

# Move the active view/selection to reflect where the edit was made
$ws.Application.Goto($ws.Range("D22"), $true)
$ws.Range("D22").Select()

# Best-effort: scroll the window so the new top-left visible cell is A19
try {
    $excel.ActiveWindow.ScrollRow = 19
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

# Best-effort: reflect the updated Excel window position/size
try {
    $excel.ActiveWindow.Left = -120
    $excel.ActiveWindow.Top = -120
    $excel.ActiveWindow.Width = 29040
    $excel.ActiveWindow.Height = 15720
} catch {
}
